$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Update crime statistics table (rows 14-29) ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("G14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("H14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("M15").Value = 200
$ws.Range("K30").Copy()
$ws.Range("M15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("N15").Value = -40

$ws.Range("C16").Value = 2

$ws.Range("E16").Value = -75

$ws.Range("G16").Value = 22

$ws.Range("H16").Value = -50

$ws.Range("I16").Value = 16

$ws.Range("J16").Value = 33

$ws.Range("K16").Value = -51.515151515151

$ws.Range("L16").Value = -5.882352941176

$ws.Range("M16").Value = 33.333333333333

$ws.Range("N16").Value = -83.673469387755

$ws.Range("C17").Value = 2

$ws.Range("D17").Value = 2
$ws.Range("J30").Copy()
$ws.Range("D17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E17").Value = 0
$ws.Range("K30").Copy()
$ws.Range("E17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("F17").Value = 12

$ws.Range("H17").Value = 0

$ws.Range("I17").Value = 21

$ws.Range("J17").Value = 17

$ws.Range("K17").Value = 23.529411764705

$ws.Range("M17").Value = 75

$ws.Range("N17").Value = -70

$ws.Range("C18").Value = 2

$ws.Range("D18").Value = 8

$ws.Range("E18").Value = -75

$ws.Range("F18").Value = 19

$ws.Range("G18").Value = 24

$ws.Range("H18").Value = -20.833333333333

$ws.Range("I18").Value = 33

$ws.Range("J18").Value = 41

$ws.Range("K18").Value = -19.512195121951

$ws.Range("L18").Value = 17.857142857142

$ws.Range("M18").Value = -21.428571428571

$ws.Range("N18").Value = -73.6

$ws.Range("C19").Value = 17

$ws.Range("D19").Value = 21

$ws.Range("E19").Value = -19.047619047619

$ws.Range("F19").Value = 66

$ws.Range("G19").Value = 75

$ws.Range("H19").Value = -12

$ws.Range("I19").Value = 98

$ws.Range("J19").Value = 102

$ws.Range("K19").Value = -3.921568627450

$ws.Range("L19").Value = 81.481481481481

$ws.Range("M19").Value = 1.030927835051

$ws.Range("N19").Value = -42.690058479532

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("D20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("E20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("G20").Value = 1

$ws.Range("N20").Value = -96.296296296296

$ws.Range("C21").Value = 23

$ws.Range("D21").Value = 39

$ws.Range("E21").Value = -41.025641025641

$ws.Range("F21").Value = 111

$ws.Range("G21").Value = 136

$ws.Range("H21").Value = -18.382352941176

$ws.Range("I21").Value = 173

$ws.Range("J21").Value = 202

$ws.Range("K21").Value = -14.356435643564

$ws.Range("L21").Value = 40.650406504065

$ws.Range("M21").Value = 2.976190476190

$ws.Range("N21").Value = -67.047619047619

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("C22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("F22").Value = 1

$ws.Range("G22").Value = 1

$ws.Range("H22").Value = 0

$ws.Range("C23").Value = 1

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("D23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("E23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("F23").Value = 6

$ws.Range("G23").Value = 11

$ws.Range("H23").Value = -45.454545454545

$ws.Range("I23").Value = 13

$ws.Range("K23").Value = -13.333333333333

$ws.Range("L23").Value = -27.777777777777

$ws.Range("M23").Value = -13.333333333333

$ws.Range("C24").Value = 23

$ws.Range("D24").Value = 68

$ws.Range("E24").Value = -66.176470588235

$ws.Range("F24").Value = 107

$ws.Range("G24").Value = 148

$ws.Range("H24").Value = -27.702702702702

$ws.Range("I24").Value = 159

$ws.Range("J24").Value = 216

$ws.Range("K24").Value = -26.388888888888

$ws.Range("L24").Value = 65.625

$ws.Range("M24").Value = -5.357142857142

$ws.Range("C25").Value = 7

$ws.Range("D25").Value = 7

$ws.Range("E25").Value = 0

$ws.Range("F25").Value = 28

$ws.Range("G25").Value = 31

$ws.Range("H25").Value = -9.677419354838

$ws.Range("I25").Value = 49

$ws.Range("J25").Value = 52

$ws.Range("K25").Value = -5.769230769230

$ws.Range("L25").Value = 113.04347826087

$ws.Range("M25").Value = 6.521739130434

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("C26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("D27").Value = 4

$ws.Range("E27").Value = -75

$ws.Range("F27").Value = 2

$ws.Range("G27").Value = 9

$ws.Range("H27").Value = -77.777777777777

$ws.Range("I27").Value = 5

$ws.Range("J27").Value = 13

$ws.Range("K27").Value = -61.538461538461

$ws.Range("L27").Value = 400

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("G28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("H28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("G29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("H29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
